# ---------------------------------------------------------------------------
# Review-response revision edits
#  1) "legacy mining facility to Tampa Bay" -> "... adjacent to Tampa Bay"
#  2) Replace the single-conference sentence (which cites the BASIS7
#     hyperlink) with a lead-in sentence followed by two new bulleted
#     citation paragraphs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Small wording fix in the quoted manuscript title --------------------
$d.Content.Find.Execute(
    "legacy mining facility to Tampa Bay",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "legacy mining facility adjacent to Tampa Bay", 2
)

# --- 2. Rework the "Preliminary results..." sentence ------------------------
# Locate paragraph 16 (the long "Response:" paragraph that currently ends
# with "... Weisberg 2022, BASIS7).")
$p16 = $d.Paragraphs.Item(16)

# Find where the sentence we are rewriting starts.
$startRng = $d.Content.Duplicate
$startRng.Find.Execute("Preliminary results were reported at a recent conference")
$sentenceStart = $startRng.Start

# The sentence (and the BASIS7 hyperlink inside it) runs to the end of the
# paragraph, just before the paragraph mark.
$sentenceEnd = $p16.Range.End - 1

$targetRng = $d.Range($sentenceStart, $sentenceEnd)
$targetRng.Text = "Preliminary results were reported at recent conferences:"

# --- 3. Insert two new bulleted citation paragraphs after paragraph 16 ------
# Create two fresh empty paragraphs right after paragraph 16.
$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertParagraphAfter()
$p17 = $d.Paragraphs.Item(17)
$p17.Range.InsertParagraphAfter()

# Fill in the first new paragraph (re-use the existing bullet list, numId 1000).
$p17 = $d.Paragraphs.Item(17)
$xml1 = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1000`"/></w:numPr></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">Weisberg, R.H. (2022). Tampa Bay Coastal Ocean Model Applications, Abstract presented at the Bay Area Scientific Information Symposium, St. Petersburg, Florida, March 2022.</w:t></w:r></w:p>"
$d.Range($p17.Range.Start, $p17.Range.End).InsertXML($xml1)

# Fill in the second new paragraph.
$p18 = $d.Paragraphs.Item(18)
$xml2 = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1000`"/></w:numPr></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">Liu, Y., Weisberg, R.H., Zheng, L., Sun, Y., Chen, J. (2021), Nowcast/Forecast of the Tampa Bay, Piney Point Effluent Plume: A Rapid Response, Abstract (OS35B-1036) presented at AGU Fall Meeting, New Orleans, Louisiana, December 2021.</w:t></w:r></w:p>"
$d.Range($p18.Range.Start, $p18.Range.End).InsertXML($xml2)

Write-Output "Edits applied."
